# PROS-6581 - CCRU - new KPI tables and POS 2019
#
# This script reproduces (as closely as the COM surface allows):
#   1. workbookView tabRatio 992 -> 993 (best-effort; engine may not persist it)
#   2. Trim the trailing space off the ten "<PREFIX>@CCH Cooler " labels
#      in column B (HM, SM, CB, CS, CNT, QSR, PTR, HRC_BTNC, HRC_CTS, HRC_RC)
#   3. Move the sheet's saved selection from B8 to B6
#   4. Slightly narrow columns A:E to match the refreshed layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Tab ratio (bookViews/workbookView@tabRatio 992 -> 993)
$excel.ActiveWindow.TabRatio = 0.993

# 2. Remove the trailing space from the "<PREFIX>@CCH Cooler " entries.
#    Each label appears in 3 consecutive rows per KPI section.
$cchCoolerCells = @(
    @("B14","B15","B16",   "HM@CCH Cooler"),
    @("B34","B35","B36",   "SM@CCH Cooler"),
    @("B54","B55","B56",   "CB@CCH Cooler"),
    @("B74","B75","B76",   "CS@CCH Cooler"),
    @("B84","B85","B86",   "CNT@CCH Cooler"),
    @("B101","B102","B103","QSR@CCH Cooler"),
    @("B118","B119","B120","PTR@CCH Cooler"),
    @("B140","B141","B142","HRC_BTNC@CCH Cooler"),
    @("B157","B158","B159","HRC_CTS@CCH Cooler"),
    @("B174","B175","B176","HRC_RC@CCH Cooler")
)

foreach ($group in $cchCoolerCells) {
    $text = $group[3]
    $ws.Range($group[0]).Value = $text
    $ws.Range($group[1]).Value = $text
    $ws.Range($group[2]).Value = $text
}

# 3. Saved selection moves from B8 to B6
[void]$ws.Range("B6").Select()

# 4. Column widths A:E are refreshed to slightly narrower values
$ws.Columns.Item(1).ColumnWidth = 41.41766666666666
$ws.Columns.Item(2).ColumnWidth = 30.21358503401357
$ws.Columns.Item(3).ColumnWidth = 26.02991156462587
$ws.Columns.Item(4).ColumnWidth = 32.106442176870765
$ws.Columns.Item(5).ColumnWidth = 38.856442176870765
